$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# B10/C10/D10 hold the "Enterprises (% of total)" percentages as text values
# (26.9 / 71.8 / 98.6) -> updated to the more precise (26.88 / 71.76 / 98.64).
# A leading apostrophe keeps the entry as text (matching the original cell
# type); re-applying the "Normal" style afterwards clears the quote-prefix
# marker Excel would otherwise add, so formatting stays exactly as before.
$ws.Range("B10").Value = "'26.88"
$ws.Range("B10").Style = "Normal"

$ws.Range("C10").Value = "'71.76"
$ws.Range("C10").Style = "Normal"

$ws.Range("D10").Value = "'98.64"
$ws.Range("D10").Style = "Normal"
